# Fill in the new "SIN1" / 332439 row on the MAIN sheet (row 40), matching
# the new data that was appended once the missing SIN information became
# available, and move the visible selection down to the newly-added row
# (mirrors the workbookView/sheetView scroll-to-bottom the author saw when
# they made this edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAIN")

$ws.Range("A40").Value = "SIN1"
$ws.Range("B40").Value = 332439

# Bring MAIN to the front and park the selection on the newly-populated
# cell, same as the saved sheetView (activeCell="A40" sqref="A40").
$ws.Select() | Out-Null
$ws.Range("A40").Select() | Out-Null
